# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Intel(R) Wi-Fi 6 AX201 160MHz - 22.160.0.4
$ws.Range("C3").Value = 4335
$ws.Range("D3").Value = 69.5

# Row 4: iwlwifi
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 6220

# Row 5: Totals
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 10555
